# Timesheet update: add 23 new rows of work entries (rows 123-145) before the
# "Totaal (uur)" total row, which shifts from row 123 down to row 146.
# Also refresh the total SUM formula and selection/scroll position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 23 new blank rows starting at row 123 (pushes the old total row,
# and everything below it, down by 23 rows).
$ws.Rows("123:145").Insert()

# Copy the formatting (styles/number formats/borders) from the last existing
# data row (now row 122) down across the newly inserted rows so they match
# the rest of the table exactly.
$ws.Range("A122:C122").Copy()
$ws.Range("A123:C145").PasteSpecial(-4122)
$ws.Range("A123:C145").RowHeight = 24.95

$ws.Range("A123").Value = 43315
$ws.Range("B123").Value = 'Voeg functionaliteit toe aan datum knoppen op home'
$ws.Range("C123").Formula = "=0.5+0.25+0.5"
$ws.Range("A124").Value = 43315
$ws.Range("B124").Value = 'Code clean up'
$ws.Range("C124").Value = 0.5
$ws.Range("A125").Value = 43316
$ws.Range("B125").Value = 'Zoeken naar nog alternatieve APIs voor events'
$ws.Range("C125").Value = 0.25
$ws.Range("A126").Value = 43316
$ws.Range("B126").Value = 'Toevoegen van Admin table en Admin Organisation pivot'
$ws.Range("C126").Value = 0.5
$ws.Range("A127").Value = 43316
$ws.Range("B127").Value = 'Organisation table aanpassen'
$ws.Range("C127").Value = 1.75
$ws.Range("A128").Value = 43319
$ws.Range("B128").Value = 'Styling van profiel aanpassen'
$ws.Range("C128").Formula = "=0.25+0.75+0.5"
$ws.Range("A129").Value = 43319
$ws.Range("B129").Value = 'Dashboard voorzien voor organisaties'
$ws.Range("C129").Formula = "=0.25+0.75+0.5+0.25"
$ws.Range("A130").Value = 43320
$ws.Range("B130").Value = 'Profile get functie aanpassen + layout'
$ws.Range("C130").Value = 1.25
$ws.Range("A131").Value = 43321
$ws.Range("B131").Value = 'Favicon instellen'
$ws.Range("C131").Value = 0.25
$ws.Range("A132").Value = 43321
$ws.Range("B132").Value = 'User toevoegen als admin bij aanmaken van organisatie'
$ws.Range("C132").Value = 1
$ws.Range("A133").Value = 43321
$ws.Range("B133").Value = 'Display admins op organisatie pagina'
$ws.Range("C133").Formula = "=1.5"
$ws.Range("A134").Value = 43321
$ws.Range("B134").Value = 'Fix bug bij favorieten'
$ws.Range("C134").Value = 0.5
$ws.Range("A135").Value = 43321
$ws.Range("B135").Value = 'Display blogposts in newsfeed op profielpagina'
$ws.Range("C135").Value = 0.5
$ws.Range("A136").Value = 43322
$ws.Range("B136").Value = 'Interests tree voorzien'
$ws.Range("C136").Value = 2.75
$ws.Range("A137").Value = 43323
$ws.Range("B137").Value = 'Default interests tonen op homepage'
$ws.Range("C137").Value = 0.75
$ws.Range("A138").Value = 43323
$ws.Range("B138").Value = 'Default interests verwerken in search parameters'
$ws.Range("C138").Value = 1.5
$ws.Range("A139").Value = 43332
$ws.Range("B139").Value = 'Kleine styling aanpassingen'
$ws.Range("C139").Value = 0.5
$ws.Range("A140").Value = 43332
$ws.Range("B140").Value = 'Minor bugfixes'
$ws.Range("C140").Formula = "=0.25+0.75"
$ws.Range("A141").Value = 43332
$ws.Range("B141").Value = 'Slug voorzien bij organisatie'
$ws.Range("C141").Value = 0.25
$ws.Range("A142").Value = 43332
$ws.Range("B142").Value = 'Interesses toevoegen aan events'
$ws.Range("C142").Value = 1
$ws.Range("A143").Value = 43332
$ws.Range("B143").Value = 'Interesses toevoegen aan users'
$ws.Range("C143").Formula = "=0.25+0.25+0.5"
$ws.Range("A144").Value = 43332
$ws.Range("B144").Value = 'Zoeken door event interesses'
$ws.Range("C144").Value = 0.5
$ws.Range("A145").Value = 43332
$ws.Range("B145").Value = 'Dossier aanpassen'
$ws.Range("C145").Value = 1.25

# Fix up the total row's SUM formula to cover the new data range (it now
# lives at row 146 after the insert above).
$ws.Range("C146").Formula = "=SUM(C2:C145)"

# Restore the sheet selection to match where the author left off editing.
$ws.Range("B142").Select()
